$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 306, pushing the existing rows 306-318 down to 307-319.
$ws.Rows.Item(306).Insert()

# The newly inserted row 306 is a duplicate (shifted copy) of what is now row 307
# (the original row 306 data), except for a handful of changed fields.
# Copy the sibling row's values across, then overwrite the changed cells.

$ws.Range("A306").Value = 7
$ws.Range("B306").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C306").Value = "Ñuble"
$ws.Range("D306").Value = 45041
$ws.Range("E306").Value = 16
$ws.Range("F306").Value = 100112024
$ws.Range("G306").Value = "Choclo"
$ws.Range("H306").Value = "Choclero"
$ws.Range("I306").Value = "Primera"
$ws.Range("J306").Value = 10000
$ws.Range("K306").Value = 400
$ws.Range("L306").Value = 400
$ws.Range("M306").Value = 400
$ws.Range("N306").Value = "$/unidad"
$ws.Range("O306").Value = "Región del Maule"
$ws.Range("P306").Value = 400
$ws.Range("Q306").Value = 1
$ws.Range("R306").Value = "Hortaliza"

# Make sure the D306 cell keeps the date/time number format used by the rest
# of column D (style index 2 in the original workbook).
$ws.Range("D306").NumberFormat = $ws.Range("D307").NumberFormat
